$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 7761.909  # H64
$ws.Cells.Item(64, 9).Value = 5127  # I64
$ws.Cells.Item(64, 10).Value = 8750  # J64
$ws.Cells.Item(64, 11).Value = 5127  # K64
$ws.Cells.Item(64, 12).Value = 8750  # L64
$ws.Cells.Item(64, 13).Value = -4879  # M64
$ws.Cells.Item(64, 14).Value = -9246  # N64

$ws.Cells.Item(67, 8).Value = 7761.909  # H67
$ws.Cells.Item(67, 9).Value = 5127  # I67
$ws.Cells.Item(67, 10).Value = 8750  # J67
$ws.Cells.Item(67, 11).Value = 5127  # K67
$ws.Cells.Item(67, 12).Value = 8750  # L67
$ws.Cells.Item(67, 13).Value = -4269  # M67
$ws.Cells.Item(67, 14).Value = -10466  # N67

$ws.Cells.Item(74, 8).Value = 8399.4  # H74
$ws.Cells.Item(74, 10).Value = 9499  # J74
$ws.Cells.Item(74, 12).Value = 9499  # L74
$ws.Cells.Item(74, 14).Value = -11371  # N74

$ws.Cells.Item(77, 8).Value = 8399.4  # H77
$ws.Cells.Item(77, 10).Value = 9499  # J77
$ws.Cells.Item(77, 12).Value = 47495  # L77
$ws.Cells.Item(77, 14).Value = -56855  # N77

$ws.Cells.Item(101, 8).Value = 478  # H101
$ws.Cells.Item(101, 9).Value = 483.6  # I101
$ws.Cells.Item(101, 10).Value = 450  # J101
$ws.Cells.Item(101, 11).Value = 1450.8  # K101
$ws.Cells.Item(101, 12).Value = 1350  # L101
$ws.Cells.Item(101, 13).Value = 171.1999999999998  # M101
$ws.Cells.Item(101, 14).Value = -4594  # N101

$ws.Cells.Item(103, 8).Value = 5174.5  # H103
$ws.Cells.Item(103, 10).Value = 5253.909  # J103
$ws.Cells.Item(103, 12).Value = 15761.727  # L103
$ws.Cells.Item(103, 14).Value = -16933.727  # N103

$ws.Cells.Item(116, 8).Value = 9410  # H116
$ws.Cells.Item(116, 9).Value = 11454  # I116
$ws.Cells.Item(116, 11).Value = 11454  # K116
$ws.Cells.Item(116, 13).Value = -8012  # M116

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 8346.333000000001  # H2
$ws.Cells.Item(2, 9).Value = 2473.5833  # I2
$ws.Cells.Item(2, 11).Value = 2473.5833  # K2
$ws.Cells.Item(2, 13).Value = -2360.5833  # M2

$ws.Cells.Item(116, 8).Value = 8346.333000000001  # H116
$ws.Cells.Item(116, 9).Value = 2473.5833  # I116
$ws.Cells.Item(116, 11).Value = 2473.5833  # K116
$ws.Cells.Item(116, 13).Value = -179.5832999999998  # M116

$ws.Cells.Item(122, 8).Value = 2563.0588  # H122
$ws.Cells.Item(122, 9).Value = 2272  # I122
$ws.Cells.Item(122, 10).Value = 4746  # J122
$ws.Cells.Item(122, 11).Value = 6816  # K122
$ws.Cells.Item(122, 12).Value = 14238  # L122
$ws.Cells.Item(122, 13).Value = -4366  # M122
$ws.Cells.Item(122, 14).Value = -19138  # N122

$ws.Cells.Item(133, 8).Value = 91840.664  # H133
$ws.Cells.Item(133, 10).Value = 91840.664  # J133
$ws.Cells.Item(133, 12).Value = 91840.664  # L133
$ws.Cells.Item(133, 14).Value = -96900.664  # N133

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 8346.333000000001  # H3
$ws.Cells.Item(3, 9).Value = 2473.5833  # I3
$ws.Cells.Item(3, 11).Value = 2473.5833  # K3
$ws.Cells.Item(3, 13).Value = -2359.5833  # M3

$ws.Cells.Item(94, 8).Value = 660.2222  # H94
$ws.Cells.Item(94, 9).Value = 422  # I94
$ws.Cells.Item(94, 10).Value = 1136.6666  # J94
$ws.Cells.Item(94, 11).Value = 422  # K94
$ws.Cells.Item(94, 12).Value = 1136.6666  # L94
$ws.Cells.Item(94, 13).Value = 29  # M94
$ws.Cells.Item(94, 14).Value = -2038.6666  # N94

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 85.69231000000001  # H7
$ws.Cells.Item(7, 9).Value = 87  # I7
$ws.Cells.Item(7, 11).Value = 87  # K7
$ws.Cells.Item(7, 13).Value = 26  # M7

$ws.Cells.Item(58, 8).Value = 3793.75  # H58
$ws.Cells.Item(58, 9).Value = 1080.1428  # I58
$ws.Cells.Item(58, 10).Value = 7592.8  # J58
$ws.Cells.Item(58, 11).Value = 1080.1428  # K58
$ws.Cells.Item(58, 12).Value = 7592.8  # L58
$ws.Cells.Item(58, 13).Value = -877.1428000000001  # M58
$ws.Cells.Item(58, 14).Value = -7998.8  # N58

$ws.Cells.Item(99, 8).Value = 4852.2856  # H99
$ws.Cells.Item(99, 9).Value = 4242.5  # I99
$ws.Cells.Item(99, 10).Value = 5665.3335  # J99
$ws.Cells.Item(99, 11).Value = 4242.5  # K99
$ws.Cells.Item(99, 12).Value = 5665.3335  # L99
$ws.Cells.Item(99, 13).Value = -2744.5  # M99
$ws.Cells.Item(99, 14).Value = -8661.333500000001  # N99

$ws.Cells.Item(122, 8).Value = 1043.5  # H122
$ws.Cells.Item(122, 9).Value = 962  # I122
$ws.Cells.Item(122, 10).Value = 1125  # J122
$ws.Cells.Item(122, 11).Value = 2886  # K122
$ws.Cells.Item(122, 12).Value = 3375  # L122
$ws.Cells.Item(122, 13).Value = -436  # M122
$ws.Cells.Item(122, 14).Value = -8275  # N122

$ws.Cells.Item(126, 8).Value = 4852.2856  # H126
$ws.Cells.Item(126, 9).Value = 4242.5  # I126
$ws.Cells.Item(126, 10).Value = 5665.3335  # J126
$ws.Cells.Item(126, 11).Value = 12727.5  # K126
$ws.Cells.Item(126, 12).Value = 16996.0005  # L126
$ws.Cells.Item(126, 13).Value = -10257.5  # M126
$ws.Cells.Item(126, 14).Value = -21936.0005  # N126

$ws.Cells.Item(136, 8).Value = 3793.75  # H136
$ws.Cells.Item(136, 9).Value = 1080.1428  # I136
$ws.Cells.Item(136, 10).Value = 7592.8  # J136
$ws.Cells.Item(136, 11).Value = 3240.4284  # K136
$ws.Cells.Item(136, 12).Value = 22778.4  # L136
$ws.Cells.Item(136, 13).Value = -690.4284000000002  # M136
$ws.Cells.Item(136, 14).Value = -27878.4  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 176.82608  # H11
$ws.Cells.Item(11, 9).Value = 166.73685  # I11
$ws.Cells.Item(11, 10).Value = 224.75  # J11
$ws.Cells.Item(11, 11).Value = 500.21055  # K11
$ws.Cells.Item(11, 12).Value = 674.25  # L11
$ws.Cells.Item(11, 13).Value = -360.21055  # M11
$ws.Cells.Item(11, 14).Value = -954.25  # N11

$ws.Cells.Item(40, 8).Value = 68.5  # H40
$ws.Cells.Item(40, 10).Value = 0  # J40
$ws.Cells.Item(40, 12).Value = 0  # L40
$ws.Cells.Item(40, 14).ClearContents()  # N40

$ws.Cells.Item(75, 8).Value = 2278.25  # H75
$ws.Cells.Item(75, 9).Value = 2500  # I75
$ws.Cells.Item(75, 10).Value = 2204.3333  # J75
$ws.Cells.Item(75, 11).Value = 7500  # K75
$ws.Cells.Item(75, 12).Value = 6612.999899999999  # L75
$ws.Cells.Item(75, 13).Value = -6502  # M75
$ws.Cells.Item(75, 14).Value = -8608.999899999999  # N75

$ws.Cells.Item(78, 8).Value = 2278.25  # H78
$ws.Cells.Item(78, 9).Value = 2500  # I78
$ws.Cells.Item(78, 10).Value = 2204.3333  # J78
$ws.Cells.Item(78, 11).Value = 22500  # K78
$ws.Cells.Item(78, 12).Value = 19838.9997  # L78
$ws.Cells.Item(78, 13).Value = -17508  # M78
$ws.Cells.Item(78, 14).Value = -29822.9997  # N78

$ws.Cells.Item(111, 8).Value = 1757.5  # H111
$ws.Cells.Item(111, 9).Value = 525  # I111
$ws.Cells.Item(111, 10).Value = 2990  # J111
$ws.Cells.Item(111, 11).Value = 1575  # K111
$ws.Cells.Item(111, 12).Value = 8970  # L111
$ws.Cells.Item(111, 13).Value = 1492  # M111
$ws.Cells.Item(111, 14).Value = -15104  # N111

$ws.Cells.Item(128, 8).Value = 609990.8  # H128
$ws.Cells.Item(128, 9).Value = 609990.8  # I128
$ws.Cells.Item(128, 11).Value = 1829972.4  # K128
$ws.Cells.Item(128, 13).Value = -1824992.4  # M128

$ws.Cells.Item(137, 8).Value = 2900  # H137
$ws.Cells.Item(137, 9).Value = 2900  # I137
$ws.Cells.Item(137, 11).Value = 8700  # K137
$ws.Cells.Item(137, 13).Value = -3600  # M137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 15335113  # H7
$ws.Cells.Item(7, 9).Value = 15714857  # I7
$ws.Cells.Item(7, 11).Value = 15714857  # K7
$ws.Cells.Item(7, 13).Value = -15714745  # M7

$ws.Cells.Item(8, 8).Value = 15335113  # H8
$ws.Cells.Item(8, 9).Value = 15714857  # I8
$ws.Cells.Item(8, 11).Value = 15714857  # K8
$ws.Cells.Item(8, 13).Value = -15714718  # M8

$ws.Cells.Item(13, 8).Value = 374  # H13
$ws.Cells.Item(13, 9).Value = 150  # I13
$ws.Cells.Item(13, 10).Value = 523.3333  # J13
$ws.Cells.Item(13, 11).Value = 150  # K13
$ws.Cells.Item(13, 12).Value = 523.3333  # L13
$ws.Cells.Item(13, 13).Value = -11  # M13
$ws.Cells.Item(13, 14).Value = -801.3333  # N13

$ws.Cells.Item(36, 8).Value = 7000  # H36
$ws.Cells.Item(36, 9).Value = 0  # I36
$ws.Cells.Item(36, 10).Value = 7000  # J36
$ws.Cells.Item(36, 11).Value = 0  # K36
$ws.Cells.Item(36, 12).Value = 7000  # L36
$ws.Cells.Item(36, 13).ClearContents()  # M36
$ws.Cells.Item(36, 14).Value = -7970  # N36

$ws.Cells.Item(80, 8).Value = 1019  # H80
$ws.Cells.Item(80, 9).Value = 1019  # I80
$ws.Cells.Item(80, 11).Value = 1019  # K80
$ws.Cells.Item(80, 13).Value = -21  # M80

$ws.Cells.Item(83, 8).Value = 1019  # H83
$ws.Cells.Item(83, 9).Value = 1019  # I83
$ws.Cells.Item(83, 11).Value = 5095  # K83
$ws.Cells.Item(83, 13).Value = -103  # M83

$ws.Cells.Item(99, 8).Value = 7017.8  # H99
$ws.Cells.Item(99, 9).Value = 4772.25  # I99
$ws.Cells.Item(99, 11).Value = 4772.25  # K99
$ws.Cells.Item(99, 13).Value = -2526.25  # M99

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 962.1111  # H22
$ws.Cells.Item(22, 9).Value = 904.2857  # I22
$ws.Cells.Item(22, 10).Value = 998.9091  # J22
$ws.Cells.Item(22, 11).Value = 904.2857  # K22
$ws.Cells.Item(22, 12).Value = 998.9091  # L22
$ws.Cells.Item(22, 13).Value = -609.2857  # M22
$ws.Cells.Item(22, 14).Value = -1588.9091  # N22

$ws.Cells.Item(27, 8).Value = 962.1111  # H27
$ws.Cells.Item(27, 9).Value = 904.2857  # I27
$ws.Cells.Item(27, 10).Value = 998.9091  # J27
$ws.Cells.Item(27, 11).Value = 904.2857  # K27
$ws.Cells.Item(27, 12).Value = 998.9091  # L27
$ws.Cells.Item(27, 13).Value = -797.2857  # M27
$ws.Cells.Item(27, 14).Value = -1212.9091  # N27

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 6230.4165  # H136
$ws.Cells.Item(136, 9).Value = 6307.8887  # I136
$ws.Cells.Item(136, 11).Value = 18923.6661  # K136
$ws.Cells.Item(136, 13).Value = -16373.6661  # M136
